$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2,3,4,8,10) have their D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado) and S (Precio $/Kg) values
# rotated between rows: row2 <-> row8, and row3 -> row4 -> row10 -> row3.
# Apply the final target values directly to each cell.

# Row 2 (was row 8's values)
$ws.Range("D2").Value = 44761
$ws.Range("M2").Value = 100
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("S2").Value = 1025

# Row 3 (was row 10's values)
$ws.Range("D3").Value = 44893
$ws.Range("M3").Value = 80
$ws.Range("N3").Value = 21000
$ws.Range("O3").Value = 22000
$ws.Range("P3").Value = 21625
$ws.Range("S3").Value = 1081

# Row 4 (was row 3's values)
$ws.Range("D4").Value = 44357
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 15000
$ws.Range("P4").Value = 14500
$ws.Range("S4").Value = 725

# Row 8 (was row 2's values)
$ws.Range("D8").Value = 44792
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 21000
$ws.Range("O8").Value = 22000
$ws.Range("P8").Value = 21500
$ws.Range("S8").Value = 1075

# Row 10 (was row 4's values)
$ws.Range("D10").Value = 44320
$ws.Range("M10").Value = 80
$ws.Range("N10").Value = 16000
$ws.Range("O10").Value = 17000
$ws.Range("P10").Value = 16500
$ws.Range("S10").Value = 825
